# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column D (municipio-nombre) moves from being described as a "measure"
# to being described as a "dimension" (refArea-style), gaining its own
# URI-Municipio mapping column; column G (depuradoras) moves the other
# way, from "dimension" to "measure", and loses its mapping-depuradoras
# reference row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: municipio-nombre -------------------------------------------------
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# --- Column G: depuradoras --------------------------------------------------------
$ws.Range("G2").Value = "iaest-measure:depuradoras"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"

# mapping-depuradoras.xlsx reference is no longer needed
$ws.Range("G5").ClearContents()
